# Auto update: 2025-12-05 18:24:43
# Applies updated scores/metrics to the 미장_비트코인_분석 sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Riot Platforms, Inc. (RIOT)
$ws.Range("K2").Value = 58.7
$ws.Range("N2").Value = 51.53902399942638

# Row 3 - Bitcoin USD (BTC-USD)
$ws.Range("D3").Value = 91442.35000000001
$ws.Range("E3").Value = 63.7
$ws.Range("F3").Value = 1.16
$ws.Range("I3").Value = 50
$ws.Range("K3").Value = 53.5
$ws.Range("N3").Value = 51.53902399942638

# Row 4 - Coinbase Global, Inc. (COIN)
$ws.Range("K4").Value = 50.5
$ws.Range("N4").Value = 51.53902399942638

# Row 5 - MARA Holdings, Inc. (MARA)
$ws.Range("K5").Value = 48.7
$ws.Range("N5").Value = 51.53902399942638

# Row 6 - Strategy Inc (MSTR)
$ws.Range("K6").Value = 34.9
$ws.Range("N6").Value = 51.53902399942638
